$d = $word.ActiveDocument

$d.Content.Find.Execute("59+10=69", $true, $false, $false, $false, $false, $true, 1, $false, "38+33=71", 2) | Out-Null
$d.Content.Find.Execute("35+44=79", $true, $false, $false, $false, $false, $true, 1, $false, "41+13=54", 2) | Out-Null
$d.Content.Find.Execute("81-36=45", $true, $false, $false, $false, $false, $true, 1, $false, "98-21=77", 2) | Out-Null
$d.Content.Find.Execute("41+12=53", $true, $false, $false, $false, $false, $true, 1, $false, "77-8=69", 2) | Out-Null
$d.Content.Find.Execute("13+29=42", $true, $false, $false, $false, $false, $true, 1, $false, "29+70=99", 2) | Out-Null
$d.Content.Find.Execute("91-17=74", $true, $false, $false, $false, $false, $true, 1, $false, "37+54=91", 2) | Out-Null
$d.Content.Find.Execute("27+67=94", $true, $false, $false, $false, $false, $true, 1, $false, "50+11=61", 2) | Out-Null
$d.Content.Find.Execute("52-9=43", $true, $false, $false, $false, $false, $true, 1, $false, "1+17=18", 2) | Out-Null
$d.Content.Find.Execute("83-65=18", $true, $false, $false, $false, $false, $true, 1, $false, "53-27=26", 2) | Out-Null
$d.Content.Find.Execute("50-6=44", $true, $false, $false, $false, $false, $true, 1, $false, "6+92=98", 2) | Out-Null
$d.Content.Find.Execute("64-50=14", $true, $false, $false, $false, $false, $true, 1, $false, "26+43=69", 2) | Out-Null
$d.Content.Find.Execute("9-2=7", $true, $false, $false, $false, $false, $true, 1, $false, "16+50=66", 2) | Out-Null
$d.Content.Find.Execute("74-21=53", $true, $false, $false, $false, $false, $true, 1, $false, "60-38=22", 2) | Out-Null
$d.Content.Find.Execute("10-6=4", $true, $false, $false, $false, $false, $true, 1, $false, "21+77=98", 2) | Out-Null
$d.Content.Find.Execute("45-13=32", $true, $false, $false, $false, $false, $true, 1, $false, "3+84=87", 2) | Out-Null
$d.Content.Find.Execute("46-2=44", $true, $false, $false, $false, $false, $true, 1, $false, "89-60=29", 2) | Out-Null
$d.Content.Find.Execute("62-3=59", $true, $false, $false, $false, $false, $true, 1, $false, "87-37=50", 2) | Out-Null
$d.Content.Find.Execute("10+60=70", $true, $false, $false, $false, $false, $true, 1, $false, "74-72=2", 2) | Out-Null
$d.Content.Find.Execute("11+15=26", $true, $false, $false, $false, $false, $true, 1, $false, "93-51=42", 2) | Out-Null
$d.Content.Find.Execute("95-12=83", $true, $false, $false, $false, $false, $true, 1, $false, "39-10=29", 2) | Out-Null
$d.Content.Find.Execute("67+5=72", $true, $false, $false, $false, $false, $true, 1, $false, "22+29=51", 2) | Out-Null
$d.Content.Find.Execute("44-8=36", $true, $false, $false, $false, $false, $true, 1, $false, "42-19=23", 2) | Out-Null
$d.Content.Find.Execute("33-8=25", $true, $false, $false, $false, $false, $true, 1, $false, "79-8=71", 2) | Out-Null
$d.Content.Find.Execute("19+6=25", $true, $false, $false, $false, $false, $true, 1, $false, "82+1=83", 2) | Out-Null
$d.Content.Find.Execute("23+14=37", $true, $false, $false, $false, $false, $true, 1, $false, "28+38=66", 2) | Out-Null
$d.Content.Find.Execute("27+36=63", $true, $false, $false, $false, $false, $true, 1, $false, "11+80=91", 2) | Out-Null
$d.Content.Find.Execute("13+84=97", $true, $false, $false, $false, $false, $true, 1, $false, "94-21=73", 2) | Out-Null
$d.Content.Find.Execute("84-55=29", $true, $false, $false, $false, $false, $true, 1, $false, "30+54=84", 2) | Out-Null
$d.Content.Find.Execute("10+9=19", $true, $false, $false, $false, $false, $true, 1, $false, "54-2=52", 2) | Out-Null
$d.Content.Find.Execute("12+85=97", $true, $false, $false, $false, $false, $true, 1, $false, "75-35=40", 2) | Out-Null
$d.Content.Find.Execute("71-60=11", $true, $false, $false, $false, $false, $true, 1, $false, "35+35=70", 2) | Out-Null
$d.Content.Find.Execute("64-27=37", $true, $false, $false, $false, $false, $true, 1, $false, "50-23=27", 2) | Out-Null
$d.Content.Find.Execute("95-67=28", $true, $false, $false, $false, $false, $true, 1, $false, "8+80=88", 2) | Out-Null
$d.Content.Find.Execute("69+17=86", $true, $false, $false, $false, $false, $true, 1, $false, "21-14=7", 2) | Out-Null
$d.Content.Find.Execute("57-51=6", $true, $false, $false, $false, $false, $true, 1, $false, "34-6=28", 2) | Out-Null
$d.Content.Find.Execute("56+3=59", $true, $false, $false, $false, $false, $true, 1, $false, "4+34=38", 2) | Out-Null
$d.Content.Find.Execute("30+9=39", $true, $false, $false, $false, $false, $true, 1, $false, "7+64=71", 2) | Out-Null
$d.Content.Find.Execute("18+17=35", $true, $false, $false, $false, $false, $true, 1, $false, "72+8=80", 2) | Out-Null
$d.Content.Find.Execute("71-14=57", $true, $false, $false, $false, $false, $true, 1, $false, "54+42=96", 2) | Out-Null
$d.Content.Find.Execute("34-7=27", $true, $false, $false, $false, $false, $true, 1, $false, "68+17=85", 2) | Out-Null
$d.Content.Find.Execute("43+45=88", $true, $false, $false, $false, $false, $true, 1, $false, "91-32=59", 2) | Out-Null
$d.Content.Find.Execute("52-10=42", $true, $false, $false, $false, $false, $true, 1, $false, "55+41=96", 2) | Out-Null
$d.Content.Find.Execute("50+0=50", $true, $false, $false, $false, $false, $true, 1, $false, "29+2=31", 2) | Out-Null
$d.Content.Find.Execute("56+24=80", $true, $false, $false, $false, $false, $true, 1, $false, "47-10=37", 2) | Out-Null
$d.Content.Find.Execute("73-4=69", $true, $false, $false, $false, $false, $true, 1, $false, "66-60=6", 2) | Out-Null
$d.Content.Find.Execute("18+21=39", $true, $false, $false, $false, $false, $true, 1, $false, "82-12=70", 2) | Out-Null
$d.Content.Find.Execute("5+42=47", $true, $false, $false, $false, $false, $true, 1, $false, "96-93=3", 2) | Out-Null
$d.Content.Find.Execute("5+64=69", $true, $false, $false, $false, $false, $true, 1, $false, "67-23=44", 2) | Out-Null
$d.Content.Find.Execute("95-48=47", $true, $false, $false, $false, $false, $true, 1, $false, "31+11=42", 2) | Out-Null
$d.Content.Find.Execute("57+18=75", $true, $false, $false, $false, $false, $true, 1, $false, "94-42=52", 2) | Out-Null
$d.Content.Find.Execute("19-3=16", $true, $false, $false, $false, $false, $true, 1, $false, "46+32=78", 2) | Out-Null
$d.Content.Find.Execute("72-47=25", $true, $false, $false, $false, $false, $true, 1, $false, "21+68=89", 2) | Out-Null
$d.Content.Find.Execute("65-57=8", $true, $false, $false, $false, $false, $true, 1, $false, "59-8=51", 2) | Out-Null
$d.Content.Find.Execute("50-25=25", $true, $false, $false, $false, $false, $true, 1, $false, "92-38=54", 2) | Out-Null
$d.Content.Find.Execute("27+53=80", $true, $false, $false, $false, $false, $true, 1, $false, "53+28=81", 2) | Out-Null
$d.Content.Find.Execute("15+3=18", $true, $false, $false, $false, $false, $true, 1, $false, "67-48=19", 2) | Out-Null
$d.Content.Find.Execute("7+61=68", $true, $false, $false, $false, $false, $true, 1, $false, "10+22=32", 2) | Out-Null
$d.Content.Find.Execute("46-4=42", $true, $false, $false, $false, $false, $true, 1, $false, "48+4=52", 2) | Out-Null
$d.Content.Find.Execute("51-7=44", $true, $false, $false, $false, $false, $true, 1, $false, "74+16=90", 2) | Out-Null
$d.Content.Find.Execute("14+69=83", $true, $false, $false, $false, $false, $true, 1, $false, "83-54=29", 2) | Out-Null
$d.Content.Find.Execute("35-28=7", $true, $false, $false, $false, $false, $true, 1, $false, "13+19=32", 2) | Out-Null
$d.Content.Find.Execute("36+5=41", $true, $false, $false, $false, $false, $true, 1, $false, "87-82=5", 2) | Out-Null
$d.Content.Find.Execute("14-9=5", $true, $false, $false, $false, $false, $true, 1, $false, "65-12=53", 2) | Out-Null
$d.Content.Find.Execute("44+9=53", $true, $false, $false, $false, $false, $true, 1, $false, "14+35=49", 2) | Out-Null
$d.Content.Find.Execute("45-17=28", $true, $false, $false, $false, $false, $true, 1, $false, "55+5=60", 2) | Out-Null
$d.Content.Find.Execute("30+61=91", $true, $false, $false, $false, $false, $true, 1, $false, "79-46=33", 2) | Out-Null
$d.Content.Find.Execute("43+31=74", $true, $false, $false, $false, $false, $true, 1, $false, "71-0=71", 2) | Out-Null
$d.Content.Find.Execute("93-47=46", $true, $false, $false, $false, $false, $true, 1, $false, "80-62=18", 2) | Out-Null
$d.Content.Find.Execute("67-46=21", $true, $false, $false, $false, $false, $true, 1, $false, "35-25=10", 2) | Out-Null
$d.Content.Find.Execute("33+64=97", $true, $false, $false, $false, $false, $true, 1, $false, "89-13=76", 2) | Out-Null
$d.Content.Find.Execute("57+25=82", $true, $false, $false, $false, $false, $true, 1, $false, "53+30=83", 2) | Out-Null
$d.Content.Find.Execute("11+88=99", $true, $false, $false, $false, $false, $true, 1, $false, "51-28=23", 2) | Out-Null
$d.Content.Find.Execute("81-6=75", $true, $false, $false, $false, $false, $true, 1, $false, "40+30=70", 2) | Out-Null
$d.Content.Find.Execute("11+68=79", $true, $false, $false, $false, $false, $true, 1, $false, "70-57=13", 2) | Out-Null
$d.Content.Find.Execute("43+6=49", $true, $false, $false, $false, $false, $true, 1, $false, "13+39=52", 2) | Out-Null
$d.Content.Find.Execute("70+6=76", $true, $false, $false, $false, $false, $true, 1, $false, "27+54=81", 2) | Out-Null
$d.Content.Find.Execute("15+47=62", $true, $false, $false, $false, $false, $true, 1, $false, "26-21=5", 2) | Out-Null
$d.Content.Find.Execute("14+40=54", $true, $false, $false, $false, $false, $true, 1, $false, "23+71=94", 2) | Out-Null
$d.Content.Find.Execute("29-13=16", $true, $false, $false, $false, $false, $true, 1, $false, "16+1=17", 2) | Out-Null
$d.Content.Find.Execute("15+10=25", $true, $false, $false, $false, $false, $true, 1, $false, "42-41=1", 2) | Out-Null
$d.Content.Find.Execute("21+71=92", $true, $false, $false, $false, $false, $true, 1, $false, "36-7=29", 2) | Out-Null
$d.Content.Find.Execute("62+5=67", $true, $false, $false, $false, $false, $true, 1, $false, "27+30=57", 2) | Out-Null
$d.Content.Find.Execute("66+26=92", $true, $false, $false, $false, $false, $true, 1, $false, "93-30=63", 2) | Out-Null
$d.Content.Find.Execute("34+21=55", $true, $false, $false, $false, $false, $true, 1, $false, "88-26=62", 2) | Out-Null
$d.Content.Find.Execute("18+65=83", $true, $false, $false, $false, $false, $true, 1, $false, "93-29=64", 2) | Out-Null
$d.Content.Find.Execute("58+36=94", $true, $false, $false, $false, $false, $true, 1, $false, "86-11=75", 2) | Out-Null
$d.Content.Find.Execute("78-8=70", $true, $false, $false, $false, $false, $true, 1, $false, "82-56=26", 2) | Out-Null
$d.Content.Find.Execute("28+52=80", $true, $false, $false, $false, $false, $true, 1, $false, "22+6=28", 2) | Out-Null
$d.Content.Find.Execute("47+29=76", $true, $false, $false, $false, $false, $true, 1, $false, "94+3=97", 2) | Out-Null
$d.Content.Find.Execute("31-3=28", $true, $false, $false, $false, $false, $true, 1, $false, "95-11=84", 2) | Out-Null
$d.Content.Find.Execute("51+25=76", $true, $false, $false, $false, $false, $true, 1, $false, "20+5=25", 2) | Out-Null
$d.Content.Find.Execute("48-23=25", $true, $false, $false, $false, $false, $true, 1, $false, "18+22=40", 2) | Out-Null
$d.Content.Find.Execute("44-33=11", $true, $false, $false, $false, $false, $true, 1, $false, "76-63=13", 2) | Out-Null
$d.Content.Find.Execute("29+31=60", $true, $false, $false, $false, $false, $true, 1, $false, "47-38=9", 2) | Out-Null
$d.Content.Find.Execute("24+67=91", $true, $false, $false, $false, $false, $true, 1, $false, "59+27=86", 2) | Out-Null
$d.Content.Find.Execute("28-21=7", $true, $false, $false, $false, $false, $true, 1, $false, "15+45=60", 2) | Out-Null
$d.Content.Find.Execute("82-27=55", $true, $false, $false, $false, $false, $true, 1, $false, "33+48=81", 2) | Out-Null
$d.Content.Find.Execute("32+9=41", $true, $false, $false, $false, $false, $true, 1, $false, "20+44=64", 2) | Out-Null
$d.Content.Find.Execute("15+79=94", $true, $false, $false, $false, $false, $true, 1, $false, "99-55=44", 2) | Out-Null
$d.Content.Find.Execute("70+19=89", $true, $false, $false, $false, $false, $true, 1, $false, "35-23=12", 2) | Out-Null
